$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "trainingaudio/01_kitipi1.wav"
$ws.Range("B2").Value = "pngimages/01_gift.png"

$ws.Range("A3").Value = "trainingaudio/16_kotapi2.wav"
$ws.Range("B3").Value = "pngimages/16_icecream.png"

$ws.Range("A4").Value = "trainingaudio/05_titopo2.wav"
$ws.Range("B4").Value = "pngimages/05_megaphone.png"

$ws.Range("A5").Value = "trainingaudio/26_kapako1.wav"
$ws.Range("B5").Value = "pngimages/26_pineapple.png"

$ws.Range("A6").Value = "trainingaudio/15_kopota3.wav"
$ws.Range("B6").Value = "pngimages/15_barrel.png"

$ws.Range("A7").Value = "trainingaudio/23_patoko1.wav"
$ws.Range("B7").Value = "pngimages/23_lemon.png"
